$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($ws, $ref, $text)
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws 'D2' '302.56'
Set-TextValue $ws 'E2' '1.95%'
Set-TextValue $ws 'D3' '44.11'
Set-TextValue $ws 'E3' '6.44%'
Set-TextValue $ws 'D4' '5.087'
Set-TextValue $ws 'E4' '0.99%'
Set-TextValue $ws 'D5' '0.07697'
Set-TextValue $ws 'B6' 'GateToken'
Set-TextValue $ws 'C6' 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
Set-TextValue $ws 'D6' '4.419'
Set-TextValue $ws 'E6' '1.49%'
Set-TextValue $ws 'B7' 'FTXToken'
Set-TextValue $ws 'C7' 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws 'D7' '1.616'
Set-TextValue $ws 'E7' '2.48%'
Set-TextValue $ws 'B8' 'MXToken'
Set-TextValue $ws 'C8' 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws 'D8' '1.047'
Set-TextValue $ws 'E8' '12.69%'
Set-TextValue $ws 'B9' 'LiechtensteinCryptoassetsExchange'
Set-TextValue $ws 'C9' 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue $ws 'D9' '0.1272'
Set-TextValue $ws 'E9' '7.50%'
Set-TextValue $ws 'B10' 'WazirX'
Set-TextValue $ws 'C10' 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue $ws 'D10' '0.1857'
Set-TextValue $ws 'E10' '2.22%'
Set-TextValue $ws 'B11' 'MandalaExchangeToken'
Set-TextValue $ws 'C11' 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue $ws 'D11' '0.09203'
Set-TextValue $ws 'E11' '3.28%'
Set-TextValue $ws 'B12' 'BitrueCoin'
Set-TextValue $ws 'C12' 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue $ws 'D12' '0.04195'
Set-TextValue $ws 'E12' '0.26%'
Set-TextValue $ws 'B13' 'BitMartToken'
Set-TextValue $ws 'C13' 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue $ws 'D13' '0.1047'
Set-TextValue $ws 'E13' '-0.35%'
Set-TextValue $ws 'B14' 'BitForexToken'
Set-TextValue $ws 'C14' 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue $ws 'D14' '0.001283'
Set-TextValue $ws 'E14' '-0.01%'
Set-TextValue $ws 'B15' 'TigerCash'
Set-TextValue $ws 'C15' 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue $ws 'D15' '0.005760'
Set-TextValue $ws 'E15' '-3.13%'
Set-TextValue $ws 'B16' 'UpBots'
Set-TextValue $ws 'C16' 'https://coinranking.com/coin/m5ozaAIK6+upbots-ubxt'
Set-TextValue $ws 'D16' '0.007489'
Set-TextValue $ws 'E16' '1,914.28%'
Set-TextValue $ws 'B17' 'LEO'
Set-TextValue $ws 'C17' 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue $ws 'D17' '3.348'
Set-TextValue $ws 'E17' '0.07%'
Set-TextValue $ws 'D19' '0.3356'
Set-TextValue $ws 'E19' '1.35%'
Set-TextValue $ws 'D20' '8.066'
Set-TextValue $ws 'E20' '2.73%'
Set-TextValue $ws 'D21' '0.1399'
Set-TextValue $ws 'E21' '-0.68%'
Set-TextValue $ws 'E22' '7.25%'
Set-TextValue $ws 'D23' '0.04197'
Set-TextValue $ws 'E23' '4.16%'
Set-TextValue $ws 'D24' '0.001282'
Set-TextValue $ws 'E24' '1.47%'
Set-TextValue $ws 'D25' '0.004411'
Set-TextValue $ws 'E25' '14.16%'
Set-TextValue $ws 'E26' '9.83%'
Set-TextValue $ws 'D38' '0.02492'
Set-TextValue $ws 'E38' '3.84%'
Set-TextValue $ws 'D39' '0.05296'
Set-TextValue $ws 'E39' '1.77%'
Set-TextValue $ws 'D40' '0.005928'
Set-TextValue $ws 'E40' '-13.76%'
Set-TextValue $ws 'D41' '0.007733'
Set-TextValue $ws 'E41' '-0.54%'
Set-TextValue $ws 'D42' '0.1350'
Set-TextValue $ws 'E42' '2.14%'
Set-TextValue $ws 'D43' '0.007346'
Set-TextValue $ws 'E43' '-0.35%'
Set-TextValue $ws 'D44' '0.007550'
Set-TextValue $ws 'E44' '5.32%'
Set-TextValue $ws 'D45' '0.3007'
Set-TextValue $ws 'E45' '-6.70%'
Set-TextValue $ws 'D46' '0.00006650'
Set-TextValue $ws 'E46' '7.23%'
Set-TextValue $ws 'E47' '0.05%'
Set-TextValue $ws 'D48' '0.04459'
Set-TextValue $ws 'E48' '-3.96%'
Set-TextValue $ws 'E49' '0.12%'
Set-TextValue $ws 'E50' '0.05%'
Set-TextValue $ws 'E51' '0.05%'
